$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 26 (shifts old row 26 -> row 27, keeping its data intact)
$ws.Rows("26:26").Insert()

# Populate new row 26 with a copy of the (old) row 25 values -
# i.e. the weekly entry that used to sit at row 25 before this week's update.
$ws.Range("A26").Value = 1
$ws.Range("B26").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C26").Value = "Arica y Parinacota"
$ws.Range("D26").Value = 44714
$ws.Range("E26").Value = 15
$ws.Range("F26").Value = 100112003
$ws.Range("G26").Value = "Ajo"
$ws.Range("H26").Value = "Chino"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 400
$ws.Range("K26").Value = 19000
$ws.Range("L26").Value = 20000
$ws.Range("M26").Value = 19500
$ws.Range("N26").Value = "$/caja 10 kilos"
$ws.Range("O26").Value = "China"
$ws.Range("P26").Value = 1950
$ws.Range("Q26").Value = 10
$ws.Range("R26").Value = "Hortaliza"

# Update row 25 with this week's new figures
$ws.Range("D25").Value = 44890
$ws.Range("K25").Value = 16000
$ws.Range("L25").Value = 17000
$ws.Range("M25").Value = 16500
$ws.Range("P25").Value = 1650
